$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the footnote that lived in column P (row 4) to column R first -
#    it must happen before the J:O -> L:Q shift below because that shift's
#    destination range overlaps column P.
# ---------------------------------------------------------------------------
$ws.Range("P4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4104) | Out-Null       # xlPasteAll
$excel.CutCopyMode = 0
$ws.Range("P4").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 2. Shift the "# Testing" block two columns to the right: J:O -> L:Q
#    (Row 1 headers keep their text, rows 2:51 keep their numeric values.)
#    Copy formats first (so previously-empty P/Q pick up the s="4" style),
#    then copy values on top.
# ---------------------------------------------------------------------------
$ws.Range("J1:O51").Copy() | Out-Null
$ws.Range("L1:Q51").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J1:O51").Copy() | Out-Null
$ws.Range("L1:Q51").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Wipe the now-stale J:K leftovers (everything else in J:O has already
#    been relocated to L:Q above), then rebuild J as the new
#    "Sum Training" column.
# ---------------------------------------------------------------------------
$ws.Range("J1:K51").Clear() | Out-Null

$ws.Range("C1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("J1").Value = "Sum Training"

for ($r = 2; $r -le 51; $r++) {
  $formula = "=SUM(C" + $r + ":H" + $r + ")"
  $ws.Range("J" + $r).Formula = $formula
}

# ---------------------------------------------------------------------------
# 4. Add the new "Sum Testing" column in S.
# ---------------------------------------------------------------------------
$ws.Range("C1").Copy() | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null         # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("S1").Value = "Sum Testing"

for ($r = 2; $r -le 51; $r++) {
  $formula = "=SUM(L" + $r + ":Q" + $r + ")"
  $ws.Range("S" + $r).Formula = $formula
}

# ---------------------------------------------------------------------------
# 5. Conditional formatting: drop the stale J2:O51 rules and recreate the
#    pair of rules on the relocated L2:Q51 range, then add the two brand new
#    "flag totals below expectation" rules on J2:J51 and S2:S51.
# ---------------------------------------------------------------------------
$ws.Range("J2:O51").FormatConditions.Delete() | Out-Null

$fc1 = $ws.Range("L2:Q51").FormatConditions.Add(8, 4, "6")
$fc1.Font.Color = 255
$fc1.Interior.Color = 11316474

$fc2 = $ws.Range("L2:Q51").FormatConditions.Add(8, 3, "6")
$fc2.Font.Color = 24832
$fc2.Interior.Color = 13561798

$fc3 = $ws.Range("J2:J51").FormatConditions.Add(8, 6, "144")
$fc3.Font.Color = 393372
$fc3.Interior.Color = 13551615

$fc4 = $ws.Range("S2:S51").FormatConditions.Add(8, 6, "36")
$fc4.Font.Color = 393372
$fc4.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# 6. View tweaks captured in the diff.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("I8").Select() | Out-Null
